# Updated symbol list on Fri Feb 10 03:44:40 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto ranking sheet. Every value in this sheet is stored as plain text
# (inline strings, no numeric formatting was ever applied), so each cell is
# forced to Text format ("@") before the new literal is written -- this stops
# Excel's normal "looks like a number/percentage" auto-conversion from
# mangling things like trailing zeros ("0.1740"), tiny magnitudes
# ("0.00006651"), or percent strings ("-4.44%"). The cell's original style is
# then restored so only the displayed text changes, matching the source data
# feed's formatting exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell known to carry the sheet's default (unstyled) look, used to restore
# each edited cell's style after the text-format trick above.
$defaultStyle = $ws.Range("D3").Style

$updates = @(
    @{ Cell = 'D2'; Value = '307.75' },
    @{ Cell = 'E2'; Value = '-4.44%' },
    @{ Cell = 'D3'; Value = '40.03' },
    @{ Cell = 'E3'; Value = '-5.54%' },
    @{ Cell = 'D4'; Value = '5.008' },
    @{ Cell = 'E4'; Value = '-4.90%' },
    @{ Cell = 'D5'; Value = '0.07678' },
    @{ Cell = 'E5'; Value = '-5.71%' },
    @{ Cell = 'D6'; Value = '4.240' },
    @{ Cell = 'E6'; Value = '-2.71%' },
    @{ Cell = 'D7'; Value = '1.624' },
    @{ Cell = 'E7'; Value = '-8.60%' },
    @{ Cell = 'D8'; Value = '0.8903' },
    @{ Cell = 'E8'; Value = '-7.15%' },
    @{ Cell = 'D9'; Value = '0.09966' },
    @{ Cell = 'E9'; Value = '-11.60%' },
    @{ Cell = 'D10'; Value = '0.1740' },
    @{ Cell = 'E10'; Value = '-6.37%' },
    @{ Cell = 'D11'; Value = '0.08942' },
    @{ Cell = 'E11'; Value = '-3.96%' },
    @{ Cell = 'D12'; Value = '0.04380' },
    @{ Cell = 'E12'; Value = '-5.56%' },
    @{ Cell = 'E13'; Value = '-0.90%' },
    @{ Cell = 'D14'; Value = '0.001272' },
    @{ Cell = 'E14'; Value = '-1.46%' },
    @{ Cell = 'D15'; Value = '0.005897' },
    @{ Cell = 'E15'; Value = '0.39%' },
    @{ Cell = 'D16'; Value = '3.357' },
    @{ Cell = 'D18'; Value = '0.3361' },
    @{ Cell = 'E18'; Value = '-0.08%' },
    @{ Cell = 'D19'; Value = '7.038' },
    @{ Cell = 'E19'; Value = '-5.81%' },
    @{ Cell = 'D20'; Value = '0.1341' },
    @{ Cell = 'E20'; Value = '-2.24%' },
    @{ Cell = 'E21'; Value = '16.53%' },
    @{ Cell = 'D22'; Value = '0.04226' },
    @{ Cell = 'E22'; Value = '0.44%' },
    @{ Cell = 'E23'; Value = '-4.99%' },
    @{ Cell = 'E24'; Value = '-5.69%' },
    @{ Cell = 'E25'; Value = '-7.10%' },
    @{ Cell = 'E26'; Value = '-1.01%' },
    @{ Cell = 'D38'; Value = '0.02360' },
    @{ Cell = 'E38'; Value = '-8.76%' },
    @{ Cell = 'D39'; Value = '0.05174' },
    @{ Cell = 'E39'; Value = '-5.20%' },
    @{ Cell = 'D40'; Value = '0.007960' },
    @{ Cell = 'E40'; Value = '1.27%' },
    @{ Cell = 'D41'; Value = '0.1323' },
    @{ Cell = 'E41'; Value = '-5.15%' },
    @{ Cell = 'D42'; Value = '0.006569' },
    @{ Cell = 'E42'; Value = '-0.96%' },
    @{ Cell = 'D43'; Value = '0.002042' },
    @{ Cell = 'E43'; Value = '-4.43%' },
    @{ Cell = 'D44'; Value = '0.007617' },
    @{ Cell = 'E44'; Value = '-12.56%' },
    @{ Cell = 'E45'; Value = '-7.21%' },
    @{ Cell = 'D46'; Value = '0.00006651' },
    @{ Cell = 'E46'; Value = '-5.63%' },
    @{ Cell = 'E47'; Value = '-0.98%' },
    @{ Cell = 'D48'; Value = '0.003768' },
    @{ Cell = 'E48'; Value = '7.60%' },
    @{ Cell = 'E49'; Value = '40.16%' },
    @{ Cell = 'E50'; Value = '-0.98%' },
    @{ Cell = 'E51'; Value = '-0.98%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = $defaultStyle
}
